# Update countries & provincias Spain
# Refreshes COVID case stats for several countries and re-ranks the
# "Azerbaiyan / Camerun / Bolivia" and "Belice / Nueva Caledonia" blocks
# to reflect their new case counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param(
        [int]$Row,
        [string]$Pais,
        [double]$Casos,
        [double]$Nuevos,
        [double]$Activos,
        [double]$Recuperados,
        [double]$Criticos,
        [double]$MuertesHoy,
        [double]$Muertes
    )

    $ws.Cells.Item($Row, 1).Value = $Pais
    $ws.Cells.Item($Row, 2).Value = $Casos
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Estados Unidos
Set-Row 4 "Estados Unidos" 1373442 5804 257773 1034668 16514 214 81001

# India
Set-Row 15 "India" 70717 3556 22549 45877 0 79 2291

# Marruecos
Set-Row 55 "Marruecos" 6281 218 2811 3282 1 0 188

# Kazajistan
Set-Row 59 "Kazajistan" 5207 117 2074 3101 33 1 32

# Azerbaiyan overtakes Camerun and Bolivia in total cases, so the three
# countries re-rank into rows 72-74.
Set-Row 72 "Azerbaiyan" 2589 70 1680 877 33 0 32
Set-Row 73 "Camerun" 2579 0 1465 1000 28 0 114
Set-Row 74 "Bolivia" 2556 119 273 2165 3 4 118

# Monaco
Set-Row 165 "Monaco" 96 0 85 7 1 0 4

# Belice overtakes Nueva Caledonia, swapping rows 192-193.
Set-Row 192 "Belice" 18 0 16 0 0 0 2
Set-Row 193 "Nueva Caledonia" 18 0 18 0 0 0 0
